# Update the build timestamp embedded in the version string throughout the workbook.
$wb = $excel.ActiveWorkbook

$oldText = "February 03 2026 17.29.55 EST"
$newText = "February 03 2026 18.05.36 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# "About" sheet: A2 (version banner) and A6 (recommended citation).
$aboutCells = @("A2", "A6")
foreach ($addr in $aboutCells) {
    $range = $wsAbout.Range($addr)
    $current = $range.Value2
    $range.Value = $current.Replace($oldText, $newText)
}

# "Boundaries and methane sources" sheet: build_version column S, rows 2-11.
for ($row = 2; $row -le 11; $row++) {
    $range = $wsData.Cells.Item($row, 19)
    $current = $range.Value2
    $range.Value = $current.Replace($oldText, $newText)
}
